# Correct error in slide 5 ("04"): shrink the content placeholder's
# height back down now that the picture beneath it no longer needs as
# much room reserved above it.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

foreach ($sh in $s.Shapes) {
    if ($sh.Id -eq 3) {
        # 3437729 EMU -> points (1 pt = 12700 EMU); nudged slightly so the
        # Single-precision COM round-trip truncates back to the exact EMU.
        $sh.Height = 270.6873474121094
    }
}
